$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 222.75
$ws.Range("I2").Value = 200
$ws.Range("J2").Value = 291
$ws.Range("K2").Value = 200
$ws.Range("L2").Value = 291
$ws.Range("M2").Value = -87
$ws.Range("N2").Value = -517
$ws.Range("H33").Value = 173.33333
$ws.Range("I33").Value = 173.33333
$ws.Range("J33").Value = 0
$ws.Range("K33").Value = 173.33333
$ws.Range("L33").Value = 0
$ws.Range("M33").Value = 55.66667000000001
$ws.Range("N33").ClearContents()
$ws.Range("H64").Value = 3172.4595
$ws.Range("I64").Value = 2950.111
$ws.Range("J64").Value = 3383.1052
$ws.Range("K64").Value = 2950.111
$ws.Range("L64").Value = 3383.1052
$ws.Range("M64").Value = -2702.111
$ws.Range("N64").Value = -3879.1052
$ws.Range("H67").Value = 3172.4595
$ws.Range("I67").Value = 2950.111
$ws.Range("J67").Value = 3383.1052
$ws.Range("K67").Value = 2950.111
$ws.Range("L67").Value = 3383.1052
$ws.Range("M67").Value = -2092.111
$ws.Range("N67").Value = -5099.1052
$ws.Range("H113").Value = 1911.96
$ws.Range("I113").Value = 1799.8334
$ws.Range("J113").Value = 1947.3684
$ws.Range("K113").Value = 1799.8334
$ws.Range("L113").Value = 1947.3684
$ws.Range("M113").Value = 1454.1666
$ws.Range("N113").Value = -8455.368399999999
$ws.Range("H114").Value = 36979.332
$ws.Range("J114").Value = 36979.332
$ws.Range("L114").Value = 36979.332
$ws.Range("N114").Value = -45657.332
$ws.Range("H115").Value = 11111545
$ws.Range("I115").Value = 11111545
$ws.Range("K115").Value = 33334635
$ws.Range("M115").Value = -33333068
$ws.Range("H127").Value = 1121.091
$ws.Range("I127").Value = 684.5
$ws.Range("J127").Value = 1370.5714
$ws.Range("K127").Value = 2053.5
$ws.Range("L127").Value = 4111.7142
$ws.Range("M127").Value = 2906.5
$ws.Range("N127").Value = -14031.7142
$ws.Range("H129").Value = 486.5
$ws.Range("I129").Value = 486.5
$ws.Range("K129").Value = 1459.5
$ws.Range("M129").Value = 3540.5
$ws.Range("H138").Value = 1564.75
$ws.Range("I138").Value = 1054.9459
$ws.Range("J138").Value = 2822.2666
$ws.Range("K138").Value = 3164.8377
$ws.Range("L138").Value = 8466.799800000001
$ws.Range("M138").Value = 1975.1623
$ws.Range("N138").Value = -18746.7998

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2261.2856
$ws.Range("I2").Value = 2225.7778
$ws.Range("J2").Value = 2325.2
$ws.Range("K2").Value = 2225.7778
$ws.Range("L2").Value = 2325.2
$ws.Range("M2").Value = -2112.7778
$ws.Range("N2").Value = -2551.2
$ws.Range("H74").Value = 4597.5938
$ws.Range("I74").Value = 4911.923
$ws.Range("J74").Value = 3235.5
$ws.Range("K74").Value = 4911.923
$ws.Range("L74").Value = 3235.5
$ws.Range("M74").Value = -4037.923
$ws.Range("N74").Value = -4983.5
$ws.Range("H77").Value = 4597.5938
$ws.Range("I77").Value = 4911.923
$ws.Range("J77").Value = 3235.5
$ws.Range("K77").Value = 24559.615
$ws.Range("L77").Value = 16177.5
$ws.Range("M77").Value = -20191.615
$ws.Range("N77").Value = -24913.5
$ws.Range("H107").Value = 25000
$ws.Range("J107").Value = 25000
$ws.Range("L107").Value = 25000
$ws.Range("N107").Value = -32680
$ws.Range("H116").Value = 2261.2856
$ws.Range("I116").Value = 2225.7778
$ws.Range("J116").Value = 2325.2
$ws.Range("K116").Value = 2225.7778
$ws.Range("L116").Value = 2325.2
$ws.Range("M116").Value = 68.22220000000016
$ws.Range("N116").Value = -6913.2
$ws.Range("H117").Value = 39700
$ws.Range("J117").Value = 39700
$ws.Range("L117").Value = 39700
$ws.Range("N117").Value = -48878
$ws.Range("H118").Value = 0
$ws.Range("J118").Value = 0
$ws.Range("L118").Value = 0
$ws.Range("N118").ClearContents()
$ws.Range("H132").Value = 1657.3729
$ws.Range("I132").Value = 1245.0444
$ws.Range("J132").Value = 2982.7144
$ws.Range("K132").Value = 3735.1332
$ws.Range("L132").Value = 8948.143199999999
$ws.Range("M132").Value = -1205.1332
$ws.Range("N132").Value = -14008.1432

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2261.2856
$ws.Range("I3").Value = 2225.7778
$ws.Range("J3").Value = 2325.2
$ws.Range("K3").Value = 2225.7778
$ws.Range("L3").Value = 2325.2
$ws.Range("M3").Value = -2111.7778
$ws.Range("N3").Value = -2553.2

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2956.0222
$ws.Range("I31").Value = 1577.3214
$ws.Range("J31").Value = 5226.8237
$ws.Range("K31").Value = 1577.3214
$ws.Range("L31").Value = 5226.8237
$ws.Range("M31").Value = -1282.3214
$ws.Range("N31").Value = -5816.8237
$ws.Range("H34").Value = 2956.0222
$ws.Range("I34").Value = 1577.3214
$ws.Range("J34").Value = 5226.8237
$ws.Range("K34").Value = 1577.3214
$ws.Range("L34").Value = 5226.8237
$ws.Range("M34").Value = -1375.3214
$ws.Range("N34").Value = -5630.8237
$ws.Range("H117").Value = 0
$ws.Range("J117").Value = 0
$ws.Range("L117").Value = 0
$ws.Range("N117").ClearContents()

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 734.76
$ws.Range("I131").Value = 461.15384
$ws.Range("J131").Value = 1031.1666
$ws.Range("K131").Value = 1383.46152
$ws.Range("L131").Value = 3093.4998
$ws.Range("M131").Value = 3656.53848
$ws.Range("N131").Value = -13173.4998

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H116").Value = 45800
$ws.Range("J116").Value = 45800
$ws.Range("L116").Value = 45800
$ws.Range("N116").Value = -54978
$ws.Range("H122").Value = 1124.0238
$ws.Range("I122").Value = 978.8919
$ws.Range("J122").Value = 2198
$ws.Range("K122").Value = 2936.6757
$ws.Range("L122").Value = 6594
$ws.Range("M122").Value = -486.6756999999998
$ws.Range("N122").Value = -11494

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2473.3845
$ws.Range("I7").Value = 2411.5557
$ws.Range("J7").Value = 2612.5
$ws.Range("K7").Value = 2411.5557
$ws.Range("L7").Value = 2612.5
$ws.Range("M7").Value = -2299.5557
$ws.Range("N7").Value = -2836.5
$ws.Range("H40").Value = 2538.077
$ws.Range("I40").Value = 2052.7632
$ws.Range("J40").Value = 3855.3572
$ws.Range("K40").Value = 2052.7632
$ws.Range("L40").Value = 3855.3572
$ws.Range("M40").Value = -1916.7632
$ws.Range("N40").Value = -4127.3572
$ws.Range("H46").Value = 1072.2106
$ws.Range("I46").Value = 811.6667
$ws.Range("J46").Value = 1192.4615
$ws.Range("K46").Value = 811.6667
$ws.Range("L46").Value = 1192.4615
$ws.Range("M46").Value = -623.6667
$ws.Range("N46").Value = -1568.4615
$ws.Range("H116").Value = 24733.334
$ws.Range("J116").Value = 24733.334
$ws.Range("L116").Value = 24733.334
$ws.Range("N116").Value = -33911.334
$ws.Range("H117").Value = 0
$ws.Range("J117").Value = 0
$ws.Range("L117").Value = 0
$ws.Range("N117").ClearContents()
$ws.Range("H126").Value = 2473.3845
$ws.Range("I126").Value = 2411.5557
$ws.Range("J126").Value = 2612.5
$ws.Range("K126").Value = 7234.6671
$ws.Range("L126").Value = 7837.5
$ws.Range("M126").Value = -4764.6671
$ws.Range("N126").Value = -12777.5
$ws.Range("H132").Value = 2070.41
$ws.Range("I132").Value = 1782.138
$ws.Range("J132").Value = 3999.6155
$ws.Range("K132").Value = 5346.414
$ws.Range("L132").Value = 11998.8465
$ws.Range("M132").Value = -2816.414
$ws.Range("N132").Value = -17058.8465

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H117").Value = 0
$ws.Range("J117").Value = 0
$ws.Range("L117").Value = 0
$ws.Range("N117").ClearContents()
$ws.Range("H122").Value = 2549.7917
$ws.Range("I122").Value = 2457.353
$ws.Range("J122").Value = 2774.2856
$ws.Range("K122").Value = 7372.059
$ws.Range("L122").Value = 8322.856800000001
$ws.Range("M122").Value = -4922.059
$ws.Range("N122").Value = -13222.8568
$ws.Range("H136").Value = 2839.75
$ws.Range("I136").Value = 2836.1914
$ws.Range("J136").Value = 2852.6155
$ws.Range("K136").Value = 8508.574200000001
$ws.Range("L136").Value = 8557.8465
$ws.Range("M136").Value = -5958.574200000001
$ws.Range("N136").Value = -13657.8465

